# Update league base data: swap the two rows in each of the two match pairs
# that had their positions corrected (rows 45/46 and rows 129/130).
# Column A (the running index) stays tied to the row; all other columns
# (B through AD) are swapped between the row pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# Columns B (2) through AD (30)
Swap-RowData 45 46 2 30
Swap-RowData 129 130 2 30
